$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Table 1 (StartNode, Relation, EndNode) - rows 3-9, columns A-C
$ws.Range("A3").Value = "PERSON1"
$ws.Range("C3").Value = "STUDY1"

$ws.Range("A4").Value = "PERSON1"
$ws.Range("C4").Value = "TREAT1"

$ws.Range("A5").Value = "STUDY1"
$ws.Range("C5").Value = "TREAT1"

$ws.Range("A6").Value = "STUDY1"
$ws.Range("C6").Value = "PROTOCOL1"

$ws.Range("A7").Value = "STUDY1"
$ws.Range("C7").Value = "TREAT2"

$ws.Range("A8").Value = "PERSON2"
$ws.Range("C8").Value = "STUDY1"

$ws.Range("A9").Value = "PERSON2"
$ws.Range("C9").Value = "TREAT2"

# Table 2 (Node, Property, Value) - rows 3-12, columns E-G
$ws.Range("E3").Value = "PERSON1"
$ws.Range("E4").Value = "PERSON1"
$ws.Range("E5").Value = "STUDY1"
$ws.Range("E6").Value = "TREAT1"
$ws.Range("E7").Value = "TREAT1"
$ws.Range("E8").Value = "STUDY1"
$ws.Range("E9").Value = "PERSON2"
$ws.Range("E10").Value = "PERSON2"
$ws.Range("E11").Value = "TREAT2"
$ws.Range("E12").Value = "PROTOCOL1"

# Update selection
$ws.Range("B15").Select()

# Apply sheet protection (password/hash will differ from original but structure matches)
$ws.Protect("password")
